# Fruta / hortaliza, semanal
# Insert a new weekly record at row 703 ("Mango" / Vega Central Mapocho de
# Santiago), pushing the existing rows 703:768 down to 704:769.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 703:768 -> 704:769, copying formatting from the row above
# (matches how the existing "Fecha" column keeps its date style).
$ws.Rows.Item(703).Insert()

$r = 703
$ws.Cells.Item($r, 1).Value = 9
$ws.Cells.Item($r, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($r, 3).Value = "Metropolitana"
$ws.Cells.Item($r, 4).Value = 45223
$ws.Cells.Item($r, 5).Value = 13
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100108
$ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($r, 9).Value = 100108002
$ws.Cells.Item($r, 10).Value = "Mango"
$ws.Cells.Item($r, 11).Value = "Sin especificar"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = 8000
$ws.Cells.Item($r, 15).Value = 9000
$ws.Cells.Item($r, 16).Value = 8400
$ws.Cells.Item($r, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item($r, 18).Value = "Brasil"
$ws.Cells.Item($r, 19).Value = 2100
$ws.Cells.Item($r, 20).Value = 4
